$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.283.25"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.910.15"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'307.91"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").Value = "'0.5352"
$ws.Range("E7").Value = "  +2.85%  "
$ws.Range("D8").Value = "'0.3824"
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").Value = "'0.07291"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'22.09"
$ws.Range("E10").Value = "  +4.27%  "
$ws.Range("D11").Value = "'0.9026"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").Value = "'0.08205"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "'95.88"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").Value = "'5.349"
$ws.Range("E14").Value = "  +1.28%  "
$ws.Range("D15").Value = "'1.002"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "'0.000008658"
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "'14.83"
$ws.Range("E17").Value = "  +1.75%  "
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").Value = "27.317.66"
$ws.Range("D20").Value = "1.152.53"
$ws.Range("E20").Value = "  -39.50%  "
$ws.Range("E21").Value = "  -0.76%  "
$ws.Range("E22").Value = "  +1.51%  "
$ws.Range("D23").Value = "'6.524"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("E24").Value = "  +2.04%  "
$ws.Range("D25").Value = "'2.288"
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").Value = "'18.27"
$ws.Range("E26").Value = "  +0.36%  "
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("D28").Value = "'117.11"
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("D29").Value = "'4.834"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'4.815"
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("D31").Value = "'0.09297"
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "'0.8391"
$ws.Range("E32").Value = "  +5.21%  "
$ws.Range("D33").Value = "'0.05066"
$ws.Range("E33").Value = "  -0.24%  "
$ws.Range("E34").Value = "  -1.11%  "
$ws.Range("D35").Value = "'3.008"
$ws.Range("E35").Value = "  +1.53%  "
$ws.Range("D36").Value = "'3.356"
$ws.Range("E36").Value = "  -1.99%  "
$ws.Range("D37").Value = "'2.697"
$ws.Range("D38").Value = "'0.5765"
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("D39").Value = "'0.02010"
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").Value = "'9.298"
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("E42").Value = "  -0.48%  "
$ws.Range("D43").Value = "'117.42"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "'0.1525"
$ws.Range("E44").Value = "  +0.49%  "
$ws.Range("D45").Value = "'0.4929"
$ws.Range("E45").Value = "  +1.49%  "
$ws.Range("E46").Value = "  +0.16%  "
$ws.Range("D47").Value = "'10.09"
$ws.Range("E47").Value = "  -0.53%  "
$ws.Range("D48").Value = "'1.639"
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("D49").Value = "'38.57"
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("D50").Value = "'0.06138"
$ws.Range("E50").Value = "  +3.12%  "
$ws.Range("D51").Value = "'63.50"
$ws.Range("E51").Value = "  -0.61%  "
